$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "Test Issue ID" column (B) and the extra 5th row so the
# sheet ends up holding a single "Link" column with 4 rows.
$ws.Range("A1:B5").ClearContents()

# Header
$ws.Range("A1").Value = "Link "

# Data rows with hyperlinks
$ws.Range("A2").Value = "https://jira.jnj.com/browse/JCVZ-1050"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://jira.jnj.com/browse/JCVZ-1050")
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("A3").Value = "https://jira.jnj.com/browse/AGQP-280"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://jira.jnj.com/browse/AGQP-280")
$ws.Range("A3").Style = "Hyperlink"

$ws.Range("A4").Value = "https://jira.jnj.com/browse/AFJX-11576"
$ws.Hyperlinks.Add($ws.Range("A4"), "https://jira.jnj.com/browse/AFJX-11576")
$ws.Range("A4").Style = "Hyperlink"

# Column width update (stored width 32.54296875 -> nearest value reachable
# through the COM ColumnWidth->stored-width rounding is 32.5)
$ws.Columns.Item(1).ColumnWidth = 31.67

# Selection change
$ws.Range("C8").Select()
